$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (columns A,C,D,E,F,G,H,K,P,Y) for rows 2-16
$data = @(
    @{A=11; C="Injuries";    D=35.033801; E=-85.13322700000001; F="2019-05-12"; G="17:38:35"; H="7900 SHALLOWFORD RD";            K="CHATTANOOGA";     P="17"; Y=5},
    @{A=12; C="Injuries";    D=35.033801; E=-85.13322700000001; F="2019-05-12"; G="17:38:01"; H="7900 SHALLOWFORD RD";            K="CHATTANOOGA";     P="17"; Y=5},
    @{A=18; C="Injuries";    D=35.002971; E=-85.284593;         F="2019-05-12"; G="15:28:45"; H="3410 6th Ave";                   K="CHATTANOOGA";     P="15"; Y=5},
    @{A=19; C="No Injuries"; D=35.002971; E=-85.284593;         F="2019-05-12"; G="15:27:57"; H="3410 6TH AVE";                   K="CHATTANOOGA";     P="15"; Y=5},
    @{A=31; C="Injuries";    D=35.029771; E=-85.253227;         F="2019-05-12"; G="14:09:19"; H="3300-3449 PINEWOOD TER";         K="CHATTANOOGA";     P="14"; Y=5},
    @{A=35; C="Injuries";    D=35.047766; E=-85.290509;         F="2019-05-12"; G="13:44:30"; H="E 3rd St / Wiehl St";            K="CHATTANOOGA";     P="13"; Y=5},
    @{A=36; C="Injuries";    D=35.047766; E=-85.290509;         F="2019-05-12"; G="13:42:56"; H="E 3rd St / Wiehl St";            K="CHATTANOOGA";     P="13"; Y=5},
    @{A=39; C="Injuries";    D=35.210828; E=-85.154957;         F="2019-05-12"; G="13:19:54"; H="9000 Dallas Hollow Rd";          K="HAMILTON COUNTY"; P="13"; Y=5},
    @{A=40; C="Injuries";    D=35.210828; E=-85.154957;         F="2019-05-12"; G="13:19:54"; H="9000 Dallas Hollow Rd";          K="HAMILTON COUNTY"; P="13"; Y=5},
    @{A=43; C="Injuries";    D=35.084813; E=-85.19765;          F="2019-05-12"; G="13:11:07"; H="4500-4519 Highway 58";          K="CHATTANOOGA";     P="13"; Y=5},
    @{A=56; C="Injuries";    D=35.080908; E=-85.204977;         F="2019-05-12"; G="10:25:18"; H="450 - 469 Highway 153 Sb";      K="CHATTANOOGA";     P="10"; Y=5},
    @{A=63; C="Injuries";    D=35.004212; E=-85.21047900000001; F="2019-05-12"; G="09:18:40"; H="140 INTERSTATE 75 NB";           K="CHATTANOOGA";     P="9";  Y=5},
    @{A=64; C="Injuries";    D=35.004212; E=-85.21047900000001; F="2019-05-12"; G="09:18:02"; H="140 INTERSTATE 75 NB";           K="CHATTANOOGA";     P="9";  Y=5},
    @{A=68; C="Entrapment";  D=35.052959; E=-85.180238;         F="2019-05-12"; G="04:30:11"; H="Noah Reid Rd / Shallowford Rd";  K="CHATTANOOGA";     P="4";  Y=5},
    @{A=69; C="Entrapment";  D=35.052959; E=-85.180238;         F="2019-05-12"; G="04:29:00"; H="7200 NOAH REID RD";             K="CHATTANOOGA";     P="4";  Y=5}
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    # Column F holds a date-shaped string ("2019-05-12"). Writing it directly
    # would make Excel infer a real date and reformat the cell, so force the
    # cell to Text first, write it, then drop back to the default style.
    $fCell = $ws.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row.F
    $fCell.NumberFormat = "General"
    $fCell.Style = "Normal"

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 11).Value = $row.K

    # Column P holds numeric-looking text ("17", "9", ...) that must stay text.
    $pCell = $ws.Cells.Item($r, 16)
    $pCell.NumberFormat = "@"
    $pCell.Value = $row.P
    $pCell.NumberFormat = "General"
    $pCell.Style = "Normal"

    $ws.Cells.Item($r, 25).Value = $row.Y

    $r++
}

# Remove the old trailing rows (17-20) that no longer exist in the new data.
$ws.Rows("17:20").Delete()
